$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update organisation contact details (Section 2: "Информация об организации")
$ws.Range("B6").Value  = "Национальный статистический комитет КР (Управление статистики домашних хозяйств) в рамках глобальной программы MICS ЮНИСЕФ"
$ws.Range("B7").Value  = "Калымбетова Ы.И."
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com"
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the selection to reflect where the editor left off
$ws.Range("B8").Select() | Out-Null
